$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Post Number column (A) ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# --- Publisher column (B), filled as a whole column first ---
$ws.Range("B2").Value = "Pox"
$ws.Range("B3").Value = "CYPost"
$ws.Range("B4").Value = "Forbees"

# --- Subject column (C), filled as a whole column next ---
$ws.Range("C2").Value = "Government"
$ws.Range("C3").Value = "Violence"
$ws.Range("C4").Value = "Health"

# --- Day column (D) ---
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3

# --- Row 2: Post 1 (Reaction, Hashtags, Headline) ---
$ws.Range("E2").Value = "Happy"
$ws.Range("F2").Value = "KeepThemIn"
$ws.Range("G2").Value = "GetsThingsDone"
$ws.Range("H2").Value = "InGoodHands"
$ws.Range("I2").Value = 75
$ws.Range("J2").Value = 16000
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = "Government Popularity All Time High"

# --- Image File Path column (M), filled as a whole column ---
$ws.Range("M2").Value = "post_001.jpg"
$ws.Range("M3").Value = "post_002.jpg"
$ws.Range("M4").Value = "post_003.jpg"

# --- Row 3: Post 2 (Reaction, Hashtags, Headline) ---
$ws.Range("E3").Value = "Sad"
$ws.Range("F3").Value = "TooSoon"
$ws.Range("G3").Value = "SoTragic"
$ws.Range("H3").Value = "PrayForThem"
$ws.Range("I3").Value = 117
$ws.Range("J3").Value = 28000
$ws.Range("K3").Value = 156
$ws.Range("L3").Value = "Child Dies From Riot Injuries"

# --- Row 4: Post 3 (Reaction, Hashtags, Headline) ---
$ws.Range("E4").Value = "Angry"
$ws.Range("F4").Value = "CrappyHealthCare"
$ws.Range("G4").Value = "OverpaidIdiots"
$ws.Range("H4").Value = "MoneyForWhat"
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 15000
$ws.Range("K4").Value = 80
$ws.Range("L4").Value = "Surgeons Walk Out"

# --- Number formats: Post Number column (A) as "000", Boost Cost column (K) as "0.00" ---
$ws.Range("A2:A4").NumberFormat = "000"
$ws.Range("K2:K4").NumberFormat = "0.00"

# --- Column width for Headline (L) ---
$ws.Columns("L").ColumnWidth = 34.7109375

# --- Sheet view: scroll position + selection ---
$ws.Range("L11").Select()
$excel.ActiveWindow.ScrollColumn = 4

# --- Page setup orientation ---
$ws.PageSetup.Orientation = 1
